# Auto-generated edit script: apply scheduled-runner price/profit updates
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 634.51
$ws.Range("I15").Value = 634.51
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1903.53
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -1734.53
$ws.Range("H39").Value = 134.38889
$ws.Range("I39").Value = 54.933334
$ws.Range("J39").Value = 531.6667
$ws.Range("K39").Value = 164.800002
$ws.Range("L39").Value = 1595.0001
$ws.Range("M39").Value = 131.199998
$ws.Range("N39").Value = -2187.0001
$ws.Range("H54").Value = 1133
$ws.Range("I54").Value = 1133
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 1133
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -647
$ws.Range("H58").Value = 953.7646999999999
$ws.Range("I58").Value = 518
$ws.Range("J58").Value = 1576.2858
$ws.Range("K58").Value = 1554
$ws.Range("L58").Value = 4728.857400000001
$ws.Range("M58").Value = -1404
$ws.Range("N58").Value = -5028.857400000001
$ws.Range("H87").Value = 29178.9
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 29178.9
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 29178.9
$ws.Range("N87").Value = -31674.9
$ws.Range("H90").Value = 29178.9
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 29178.9
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 87536.70000000001
$ws.Range("N90").Value = -100016.7
$ws.Range("H100").Value = 14286431
$ws.Range("I100").Value = 17241902
$ws.Range("J100").Value = 1650
$ws.Range("K100").Value = 17241902
$ws.Range("L100").Value = 1650
$ws.Range("M100").Value = -17241361
$ws.Range("N100").Value = -2732
$ws.Range("H106").Value = 8992.6875
$ws.Range("I106").Value = 10990.25
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 10990.25
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -10359.25
$ws.Range("N106").Value = -4262

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1087.1052
$ws.Range("I2").Value = 803.6667
$ws.Range("J2").Value = 2150
$ws.Range("K2").Value = 803.6667
$ws.Range("L2").Value = 2150
$ws.Range("M2").Value = -690.6667
$ws.Range("N2").Value = -2376
$ws.Range("H116").Value = 1087.1052
$ws.Range("I116").Value = 803.6667
$ws.Range("J116").Value = 2150
$ws.Range("K116").Value = 803.6667
$ws.Range("L116").Value = 2150
$ws.Range("M116").Value = 1490.3333
$ws.Range("N116").Value = -6738
$ws.Range("H132").Value = 2568.5
$ws.Range("I132").Value = 1917.2572
$ws.Range("J132").Value = 3768.158
$ws.Range("K132").Value = 5751.7716
$ws.Range("L132").Value = 11304.474
$ws.Range("M132").Value = -3221.7716
$ws.Range("N132").Value = -16364.474

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1087.1052
$ws.Range("I3").Value = 803.6667
$ws.Range("J3").Value = 2150
$ws.Range("K3").Value = 803.6667
$ws.Range("L3").Value = 2150
$ws.Range("M3").Value = -689.6667
$ws.Range("N3").Value = -2378
$ws.Range("H94").Value = 25000650
$ws.Range("I94").Value = 25000650
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 25000650
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -25000199
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4851.9688
$ws.Range("I58").Value = 1080.5
$ws.Range("J58").Value = 13149.2
$ws.Range("K58").Value = 1080.5
$ws.Range("L58").Value = 13149.2
$ws.Range("M58").Value = -877.5
$ws.Range("N58").Value = -13555.2
$ws.Range("H74").Value = 33000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 33000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 33000
$ws.Range("N74").Value = -34748
$ws.Range("H77").Value = 33000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 33000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 99000
$ws.Range("N77").Value = -107736
$ws.Range("H108").Value = 31115
$ws.Range("I108").Value = 20621
$ws.Range("J108").Value = 33738.5
$ws.Range("K108").Value = 20621
$ws.Range("L108").Value = 33738.5
$ws.Range("M108").Value = -16781
$ws.Range("N108").Value = -41418.5
$ws.Range("H132").Value = 1568.1351
$ws.Range("I132").Value = 1245.2963
$ws.Range("J132").Value = 2439.8
$ws.Range("K132").Value = 3735.8889
$ws.Range("L132").Value = 7319.400000000001
$ws.Range("M132").Value = -1205.8889
$ws.Range("N132").Value = -12379.4
$ws.Range("H135").Value = 35548.75
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 35548.75
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 35548.75
$ws.Range("N135").Value = -45688.75
$ws.Range("H136").Value = 4851.9688
$ws.Range("I136").Value = 1080.5
$ws.Range("J136").Value = 13149.2
$ws.Range("K136").Value = 3241.5
$ws.Range("L136").Value = 39447.60000000001
$ws.Range("M136").Value = -691.5
$ws.Range("N136").Value = -44547.60000000001
$ws.Range("H141").Value = 618791
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 618791
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 618791
$ws.Range("N141").Value = -629151

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 733
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 733
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 2199
$ws.Range("N127").Value = -12119

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20458368
$ws.Range("I70").Value = 14709811
$ws.Range("J70").Value = 40003460
$ws.Range("K70").Value = 14709811
$ws.Range("L70").Value = 40003460
$ws.Range("M70").Value = -14709541
$ws.Range("N70").Value = -40004000
$ws.Range("H73").Value = 20458368
$ws.Range("I73").Value = 14709811
$ws.Range("J73").Value = 40003460
$ws.Range("K73").Value = 14709811
$ws.Range("L73").Value = 40003460
$ws.Range("M73").Value = -14708875
$ws.Range("N73").Value = -40005332
$ws.Range("H102").Value = 1542.9395
$ws.Range("I102").Value = 1539.0454
$ws.Range("J102").Value = 1550.7273
$ws.Range("K102").Value = 1539.0454
$ws.Range("L102").Value = 1550.7273
$ws.Range("M102").Value = 82.95460000000003
$ws.Range("N102").Value = -4794.7273
$ws.Range("H126").Value = 1974.2609
$ws.Range("I126").Value = 1807.6666
$ws.Range("J126").Value = 2286.625
$ws.Range("K126").Value = 5422.9998
$ws.Range("L126").Value = 6859.875
$ws.Range("M126").Value = -2952.9998
$ws.Range("N126").Value = -11799.875
$ws.Range("H132").Value = 4151.421
$ws.Range("I132").Value = 4407.1
$ws.Range("J132").Value = 3867.3333
$ws.Range("K132").Value = 13221.3
$ws.Range("L132").Value = 11601.9999
$ws.Range("M132").Value = -10691.3
$ws.Range("N132").Value = -16661.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1138.091
$ws.Range("I61").Value = 1001.2222
$ws.Range("J61").Value = 1754
$ws.Range("K61").Value = 1001.2222
$ws.Range("L61").Value = 1754
$ws.Range("M61").Value = -799.2222
$ws.Range("N61").Value = -2158
$ws.Range("H68").Value = 1186.8636
$ws.Range("I68").Value = 1186.8636
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1186.8636
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -437.8635999999999
$ws.Range("H71").Value = 1186.8636
$ws.Range("I71").Value = 1186.8636
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 5934.317999999999
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -2190.317999999999
$ws.Range("H113").Value = 1138.091
$ws.Range("I113").Value = 1001.2222
$ws.Range("J113").Value = 1754
$ws.Range("K113").Value = 1001.2222
$ws.Range("L113").Value = 1754
$ws.Range("M113").Value = 1168.7778
$ws.Range("N113").Value = -6094
$ws.Range("H122").Value = 31251652
$ws.Range("I122").Value = 41668200
$ws.Range("J122").Value = 2002.5
$ws.Range("K122").Value = 125004600
$ws.Range("L122").Value = 6007.5
$ws.Range("M122").Value = -125002150
$ws.Range("N122").Value = -10907.5
$ws.Range("H132").Value = 2824.9565
$ws.Range("I132").Value = 2644.3076
$ws.Range("J132").Value = 3059.8
$ws.Range("K132").Value = 7932.9228
$ws.Range("L132").Value = 9179.400000000001
$ws.Range("M132").Value = -5402.9228
$ws.Range("N132").Value = -14239.4
$ws.Range("H136").Value = 1983.6072
$ws.Range("I136").Value = 1787.1904
$ws.Range("J136").Value = 2572.8572
$ws.Range("K136").Value = 5361.5712
$ws.Range("L136").Value = 7718.571599999999
$ws.Range("M136").Value = -2811.5712
$ws.Range("N136").Value = -12818.5716
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 585.74194
$ws.Range("I107").Value = 346.07144
$ws.Range("J107").Value = 783.1177
$ws.Range("K107").Value = 1038.21432
$ws.Range("L107").Value = 2349.3531
$ws.Range("M107").Value = 881.78568
$ws.Range("N107").Value = -6189.3531
$ws.Range("H122").Value = 33785388
$ws.Range("I122").Value = 37880384
$ws.Range("J122").Value = 1697.5
$ws.Range("K122").Value = 113641152
$ws.Range("L122").Value = 5092.5
$ws.Range("M122").Value = -113638702
$ws.Range("N122").Value = -9992.5
$ws.Range("H126").Value = 83334250
$ws.Range("I126").Value = 90909910
$ws.Range("J126").Value = 2005
$ws.Range("K126").Value = 272729730
$ws.Range("L126").Value = 6015
$ws.Range("M126").Value = -272727260
$ws.Range("N126").Value = -10955
$ws.Range("H132").Value = 2476.1943
$ws.Range("I132").Value = 2221.3462
$ws.Range("J132").Value = 3138.8
$ws.Range("K132").Value = 6664.0386
$ws.Range("L132").Value = 9416.400000000001
$ws.Range("M132").Value = -4134.0386
$ws.Range("N132").Value = -14476.4
$ws.Range("H136").Value = 1512.6
$ws.Range("I136").Value = 1058.6154
$ws.Range("J136").Value = 2355.7144
$ws.Range("K136").Value = 3175.8462
$ws.Range("L136").Value = 7067.1432
$ws.Range("M136").Value = -625.8462
$ws.Range("N136").Value = -12167.1432
